$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Step 1: stash a style-donor cell for each existing style into a scratch area (rows 201+) ---
$ws.Range("B2").Copy($ws.Range("B201"))
$ws.Range("B20").Copy($ws.Range("B202"))
$ws.Range("B21").Copy($ws.Range("B203"))
$ws.Range("B18").Copy($ws.Range("B204"))
$ws.Range("B12").Copy($ws.Range("B205"))
$ws.Range("B28").Copy($ws.Range("B206"))
$ws.Range("D29").Copy($ws.Range("B207"))
$ws.Range("D31").Copy($ws.Range("B208"))
$ws.Range("D19").Copy($ws.Range("B209"))
$ws.Range("E19").Copy($ws.Range("B210"))
$ws.Range("E21").Copy($ws.Range("B211"))

# --- Step 2: synthesize the 2 brand-new styles (center aligned, bordered) from donors ---
$ws.Range("B203").Copy($ws.Range("B212"))
$ws.Range("B212").HorizontalAlignment = -4108
$ws.Range("B211").Copy($ws.Range("B213"))
$ws.Range("B213").HorizontalAlignment = -4108

# --- Step 3: clear the existing used range completely ---
$ws.Range("A1:F40").Clear()

# --- Step 4: write every destination cell: copy style from scratch donor (if any), then set value/formula ---
$ws.Range("B201").Copy($ws.Range("B2"))
$ws.Range("B2").Value = "OpenTBS demo"
$ws.Range("B4").Value = "Hello [onshow.yourname],"
$ws.Range("B6").Value = "This is a demo of the OpenTBS plugin."
$ws.Range("B7").Value = "The current document has been generated at [onshow..now;frm='yyyy-mm-dd hh:nn:ss']"
$ws.Range("B8").Value = "PHP version : [onshow..cst.PHP_VERSION]"
$ws.Range("B9").Value = "TBS version : [onshow..version]"
$ws.Range("B205").Copy($ws.Range("B11"))
$ws.Range("B11").Value = "You may consider the following before building your own Microsoft Excel template:"
$ws.Range("B205").Copy($ws.Range("B12"))
$ws.Range("B12").Value = "* Do not use a formula in a cell that may have its position changed after the merge (for example under a TBS block). Otherwise Excel will raise an error message."
$ws.Range("B205").Copy($ws.Range("B13"))
$ws.Range("B13").Value = "    This is because the location of formulas are saved a second time in another sub-file for the order of evaluation."
$ws.Range("B205").Copy($ws.Range("B14"))
$ws.Range("B14").Value = "* If a formula uses a reference to a cell that has moved during the merge, then the reference will not be arraged to be the new cell reference. "
$ws.Range("B205").Copy($ws.Range("B15"))
$ws.Range("B15").Value = "* You cannot change picture using `"ope=changepic`". This is because drawing information are not saved directly in the sheet."
$ws.Range("B204").Copy($ws.Range("B17"))
$ws.Range("B17").Value = "Example #1: merging data with rows"
$ws.Range("B209").Copy($ws.Range("D18"))
$ws.Range("D18").Value = "Total:"
$ws.Range("B210").Copy($ws.Range("E18"))
$ws.Range("E18").Formula = "=SUM(E20:E2004)"
$ws.Range("B202").Copy($ws.Range("B19"))
$ws.Range("B19").Value = "First Name"
$ws.Range("B202").Copy($ws.Range("C19"))
$ws.Range("C19").Value = "Name"
$ws.Range("B202").Copy($ws.Range("D19"))
$ws.Range("D19").Value = "Membership number"
$ws.Range("B202").Copy($ws.Range("E19"))
$ws.Range("E19").Value = "Score"
$ws.Range("B203").Copy($ws.Range("B20"))
$ws.Range("B20").Value = "[a.firstname;block=row]"
$ws.Range("B203").Copy($ws.Range("C20"))
$ws.Range("C20").Value = "[a.name]"
$ws.Range("B203").Copy($ws.Range("D20"))
$ws.Range("D20").Value = "[a.number]"
$ws.Range("B211").Copy($ws.Range("E20"))
$ws.Range("E20").Value = "[a.score;ope=xlsxNum]"
$ws.Range("B204").Copy($ws.Range("B22"))
$ws.Range("B22").Value = "Example #2: merging data with columns"
$ws.Range("B202").Copy($ws.Range("B24"))
$ws.Range("B24").Value = "First Name:"
$ws.Range("B212").Copy($ws.Range("C24"))
$ws.Range("C24").Value = "[b1.firstname;block=c]"
$ws.Range("B202").Copy($ws.Range("B25"))
$ws.Range("B25").Value = "Score"
$ws.Range("B213").Copy($ws.Range("C25"))
$ws.Range("C25").Value = "[b2.score;block=c;ope=xlsxNum]"
$ws.Range("B204").Copy($ws.Range("B27"))
$ws.Range("B27").Value = "Example #3: change the type data in a cell"
$ws.Range("B205").Copy($ws.Range("B29"))
$ws.Range("B29").Value = "When you put a TBS field into a cell, then by default Excel assumes the cell has a string content and will not use the format you expect for the cell."
$ws.Range("B205").Copy($ws.Range("B30"))
$ws.Range("B30").Value = "But you can change the type of data in a cell using parameter « ope ». Supported types are listed in the examples below."
$ws.Range("B206").Copy($ws.Range("B32"))
$ws.Range("B32").Value = "Type of data"
$ws.Range("B206").Copy($ws.Range("C32"))
$ws.Range("C32").Value = "Parameter ope"
$ws.Range("B206").Copy($ws.Range("D32"))
$ws.Range("D32").Value = "Example"
$ws.Range("B33").Value = "Number"
$ws.Range("C33").Value = "xlsxNum"
$ws.Range("B207").Copy($ws.Range("D33"))
$ws.Range("D33").Value = "[onshow.x_num;ope=xlsxNum]"
$ws.Range("B34").Value = "Boolean"
$ws.Range("C34").Value = "xlsxBool"
$ws.Range("D34").Value = "[onshow.x_bt;ope=xlsxBool]"
$ws.Range("B35").Value = "Date/time"
$ws.Range("C35").Value = "xlsxDate"
$ws.Range("B208").Copy($ws.Range("D35"))
$ws.Range("D35").Value = "[onshow.x_dt;ope=xlsxDate]"

# --- Step 5: clear scratch area ---
$ws.Range("A200:F220").Clear()

# --- Step 6: row heights for the two thin spacer rows ---
$ws.Rows.Item(3).RowHeight = 8.25
$ws.Rows.Item(5).RowHeight = 7.5

# --- Step 7: selection ---
$ws.Range("B4").Select()
